$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("ferramenta que de respaldo as tomadas", $true, $false, $false, $false, $false, `
              $true, 1, $false, "ferramenta que respalde as tomadas", 2)
